# Applies the cryptos list update (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds plain-text figures (not numbers); keep them as text
# so values like '3.40' or '0.0000359' are preserved exactly, as in the source diff.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '68.203.67'
$ws.Range("E2").Value = '  +1.78%  '
$ws.Range("D3").Value = '3.924.33'
$ws.Range("E3").Value = '  +0.82%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '488.37'
$ws.Range("E5").Value = '  +4.27%  '
$ws.Range("D6").Value = '146.59'
$ws.Range("E6").Value = '  +2.34%  '
$ws.Range("D7").Value = '0.627'
$ws.Range("E7").Value = '  +0.38%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").Value = '0.731'
$ws.Range("E9").Value = '  -0.68%  '
$ws.Range("D10").Value = '0.169'
$ws.Range("E10").Value = '  +3.65%  '
$ws.Range("D11").Value = '0.0000359'
$ws.Range("E11").Value = '  +7.43%  '
$ws.Range("D12").Value = '42.85'
$ws.Range("E12").Value = '  +0.09%  '
$ws.Range("E13").Value = '  +3.53%  '
$ws.Range("D14").Value = '4.554.87'
$ws.Range("E14").Value = '  +1.02%  '
$ws.Range("D15").Value = '14.90'
$ws.Range("E15").Value = '  -1.97%  '
$ws.Range("D16").Value = '3.910.82'
$ws.Range("E16").Value = '  +1.06%  '
$ws.Range("E17").Value = '  -0.13%  '
$ws.Range("D18").Value = '20.12'
$ws.Range("E18").Value = '  +1.20%  '
$ws.Range("E19").Value = '  -1.52%  '
$ws.Range("D20").Value = '68.376.91'
$ws.Range("E20").Value = '  +1.81%  '
$ws.Range("D21").Value = '446.63'
$ws.Range("E21").Value = '  +3.70%  '
$ws.Range("E22").Value = '  +0.72%  '
$ws.Range("D23").Value = '3.40'
$ws.Range("E23").Value = '  +1.79%  '
$ws.Range("D24").Value = '88.64'
$ws.Range("E24").Value = '  +0.41%  '
$ws.Range("D25").Value = '11.62'
$ws.Range("E25").Value = '  +15.87%  '
$ws.Range("D26").Value = '10.95'
$ws.Range("E26").Value = '  +14.13%  '
$ws.Range("D27").Value = '3.64'
$ws.Range("E27").Value = '  +3.31%  '
$ws.Range("D28").Value = '39.11'
$ws.Range("E28").Value = '  +0.62%  '
$ws.Range("E29").Value = '  +2.41%  '
$ws.Range("B30").Value = 'Bittensor'
$ws.Range("C30").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D30").Value = '697.10'
$ws.Range("E30").Value = '  -4.49%  '
$ws.Range("B31").Value = 'Cosmos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D31").Value = '13.51'
$ws.Range("E31").Value = '  -1.35%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = '0.131'
$ws.Range("E32").Value = '  +0.48%  '
$ws.Range("E33").Value = '  +2.74%  '
$ws.Range("D34").Value = '0.0₃0939'
$ws.Range("E34").Value = '  +21.51%  '
$ws.Range("D35").Value = '41.84'
$ws.Range("E35").Value = '  -2.45%  '
$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").Value = '5.82'
$ws.Range("E36").Value = '  +8.46%  '
$ws.Range("B37").Value = 'OKB'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D37").Value = '58.92'
$ws.Range("E37").Value = '  +1.77%  '
$ws.Range("E38").Value = '  -4.09%  '
$ws.Range("E39").Value = '  +0.06%  '
$ws.Range("D40").Value = '0.0479'
$ws.Range("E40").Value = '  +1.00%  '
$ws.Range("D41").Value = '0.375'
$ws.Range("E41").Value = '  +11.11%  '
$ws.Range("D42").Value = '2.83'
$ws.Range("E42").Value = '  +13.04%  '
$ws.Range("E43").Value = '  -3.33%  '
$ws.Range("E44").Value = '  +5.87%  '
$ws.Range("E45").Value = '  +2.23%  '
$ws.Range("E46").Value = '  -0.07%  '
$ws.Range("E47").Value = '  +1.31%  '
$ws.Range("E48").Value = '  -0.54%  '
$ws.Range("D49").Value = '146.26'
$ws.Range("E49").Value = '  +1.88%  '
$ws.Range("E50").Value = '  +0.12%  '
$ws.Range("E51").Value = '  -1.27%  '

Write-Output "Updated cryptos list"
